$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) values for rows 2-49 based on latest crypto data
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.209.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0628"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.927.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.704.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.228.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "237.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0744"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.37%  "
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.551.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.950"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.603"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  -0.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.836.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  +5.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.30%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.19%  "
